$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.1290141311112389
$ws.Range("B2").Value = 0.2163186755499787

$ws.Range("A3").Value = 0.008458807205688557
$ws.Range("B3").Value = 0.1670046350658946

$ws.Range("A4").Value = 0.00072989691808817
$ws.Range("B4").Value = -0.1114220292785898

$ws.Range("A5").Value = 0.04771530505415983
$ws.Range("B5").Value = 0.1354491558372081

$ws.Range("A6").Value = -0.1191461108926047
$ws.Range("B6").Value = -0.04532125985197404

$ws.Range("A7").Value = 0.584321372854409
$ws.Range("B7").Value = 0.4812275335469187

$ws.Range("A8").Value = 0.2297604868564154
$ws.Range("B8").Value = -0.03626881835560345
